# Generate Report for Handoff
# b.md has been handed off again (new handoff xliff files generated for
# both zh-cn and de-de) and the existing handback is now stale relative
# to it, so we record the new "Ready for handoff" status, the new
# handoff file names / timestamps, and an error detail noting the
# handback file is out of date.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d3044dc8dce3f6efc0da6719d5c996541ea3493/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8fd676338d7b6e161ff81be4bd5cb6c15e1dbecb/e2e/b.md."

# --- Overview sheet: update the status/date summary columns for b.md ---
$wsOverview = $wb.Worksheets("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 04:36:39"

# --- zh-cn sheet: b.md (row 3) gets a fresh handoff ---
$wsZh = $wb.Worksheets("zh-cn")
$wsZh.Range("P1").ColumnWidth = 39.15
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-18 04:36:34"
$wsZh.Range("P3").Value = $errorDetail

# --- de-de sheet: b.md (row 3) gets a fresh handoff ---
$wsDe = $wb.Worksheets("de-de")
$wsDe.Range("P1").ColumnWidth = 39.15
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-18 04:36:39"
$wsDe.Range("P3").Value = $errorDetail
